# Apply edit to slide 15: add body text + reposition content placeholder,
# and reposition/resize the picture.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)

# EMU -> points conversion (914400 EMU per inch, 72 points per inch)
$emuPerPoint = 12700

# --- Shape 2: "Tartalom helye 2" content placeholder ---
$body = $s.Shapes.Item(2)

$body.TextFrame.TextRange.Text = "A SOAP egy XML-alapú üzenetküldési protokoll, amely lassúnak és összetettnek bizonyult, ezért fejlesztették ki az egyszerűbb, XML-t nem igénylő REST API-t, amely HTTP-t használva könnyen skálázható, gyors és népszerű webszolgáltatásokhoz; az RPC, XML-RPC és JSON-RPC ehhez hasonlóan távoli eljáráshívásra szolgál, de különböznek formázásuk és rugalmasságuk tekintetében.`r" 

$body.Left = 308219 / $emuPerPoint
$body.Top = 205954 / $emuPerPoint
$body.Width = 8596668 / $emuPerPoint
$body.Height = 3880773 / $emuPerPoint

# --- Shape 3: picture "Kép 3" ---
$pic = $s.Shapes.Item(3)

$pic.Left = 402670 / $emuPerPoint
$pic.Top = 1930400 / $emuPerPoint
$pic.Width = 9910199 / $emuPerPoint
$pic.Height = 4930116 / $emuPerPoint
